$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# Ensure G:H columns use the time (h:mm) number format to match existing rows
$ws.Range("G24:H28").NumberFormat = "h:mm"

# Row 24
$ws.Range("E24").Value = "Reviewer"
$ws.Range("F24").Value = "Testsuite OC1"
$ws.Range("G24").Value = 0.44791666666666669
$ws.Range("H24").Value = 0.46875
$ws.Range("I24").Value = "30min"

# Row 25
$ws.Range("E25").Value = "Reviewer"
$ws.Range("F25").Value = "SD for OC1"
$ws.Range("G25").Value = 0.46875
$ws.Range("H25").Value = 0.5
$ws.Range("I25").Value = "45min"

# Row 26
$ws.Range("E26").Value = "Designer"
$ws.Range("F26").Value = "SD fir OC1"
$ws.Range("G26").Value = 0.5
$ws.Range("H26").Value = 0.55208333333333337
$ws.Range("I26").Value = "1t 15m"

# Row 27
$ws.Range("E27").Value = "Reviewer"
$ws.Range("F27").Value = "DCD for OC1"
$ws.Range("G27").Value = 0.55208333333333337
$ws.Range("H27").Value = 0.57986111111111105
$ws.Range("I27").Value = "40m"

# Row 28
$ws.Range("E28").Value = "Reviewer"
$ws.Range("F28").Value = "Testsuite OC2"
$ws.Range("G28").Value = 0.57986111111111105
$ws.Range("H28").Value = 0.59375
$ws.Range("I28").Value = "20m"

# Restore the active-cell selection seen in the final sheet state
[void]$ws.Range("F14").Select()
